# Add institution validity start/end date rows (rows 23-26) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: institution / dfi / start_date / 2012
$ws.Range("A23").Value = 1747340787
$ws.Range("B23").Value = "update"
$ws.Range("C23").Value = "institution"
$ws.Range("D23").Value = "dfi"
$ws.Range("F23").Value = "start_date"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "2012"
$ws.Range("H23").ClearFormats()

# Row 24: institution / dfi / end_date / 2023
$ws.Range("A24").Value = 1747340787
$ws.Range("B24").Value = "update"
$ws.Range("C24").Value = "institution"
$ws.Range("D24").Value = "dfi"
$ws.Range("F24").Value = "end_date"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "2023"
$ws.Range("H24").ClearFormats()

# Row 25: institution / dff / start_date / 2010/10
$ws.Range("A25").Value = 1747340787
$ws.Range("B25").Value = "update"
$ws.Range("C25").Value = "institution"
$ws.Range("D25").Value = "dff"
$ws.Range("F25").Value = "start_date"
$ws.Range("H25").Value = "2010/10"

# Row 26: institution / seco / end_date / 2021/04
$ws.Range("A26").Value = 1747340787
$ws.Range("B26").Value = "update"
$ws.Range("C26").Value = "institution"
$ws.Range("D26").Value = "seco"
$ws.Range("F26").Value = "end_date"
$ws.Range("H26").Value = "2021/04"
